# "correction and print result"
#
# The sheet lists production types with their recipe's output item/quantity.
# Row 2 ("制造" / Assembling) incorrectly listed "白糖" (white sugar) x600 as
# its produced item; correct it to "电磁涡轮" (electromagnetic turbine) with
# a quantity of 0.1. Row 3 ("冶炼" / Smelting) never actually produces an
# item/quantity pair, so its stray "白糖" x600 values are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the produced item + quantity for row 2 (Assembling).
$ws.Range("C2").Value = "电磁涡轮"
$ws.Range("D2").Value = 0.1

# Row 3 (Smelting) shouldn't have had an item/quantity at all - clear it.
$ws.Range("C3:D3").ClearContents()

# Leave the selection on the corrected quantity cell, matching the saved file.
$ws.Range("D3").Select()
